$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "If this command is executed with an *attribute* that is
# already..." -> "This command toggles the state of an *attribute* within
# the set: thus if this command is executed with an *attribute* that is
# already in the excluded set, it is removed from the set, otherwise it is
# added. By default, ..."
# ---------------------------------------------------------------------------

# 1a. Rewrite the lead-in sentence.
$r1 = $d.Content
$r1.Find.Execute("If this command is executed with an ", $true, $false, $false, $false, $false, $true, 1, $false, "This command toggles the state of an ", 2) | Out-Null

# 1b. Rewrite the tail (after the first italic "attribute") and leave a
# placeholder marker for the second italic "attribute" run.
$r2 = $d.Content
$r2.Find.Execute(" that is already in the excluded set, it is removed from the set. By default, ", $true, $false, $false, $false, $false, $true, 1, $false, " within the set: thus if this command is executed with an <<ATTR>> that is already in the excluded set, it is removed from the set, otherwise it is added. By default, ", 2) | Out-Null

# 1c. Replace the placeholder with an italicized "attribute" run.
$r3 = $d.Content
$r3.Find.Execute("<<ATTR>>", $true, $false, $false, $false, $false, $true, 1, $false, "attribute", 2) | Out-Null
$r3.Font.Italic = 1

# ---------------------------------------------------------------------------
# Change 2: move the "_GoBack" bookmark from just before the
# "(E1 ^command C1 ^present-id 4 ^result R2)" code line to right after the
# "Or, visually:" paragraph. Since bookmark names are unique, re-adding a
# bookmark named "_GoBack" moves it (removing it from its old spot).
# ---------------------------------------------------------------------------

$rb = $d.Content
$rb.Find.Execute("Or, visually:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rb.Collapse(0)
# Use a temporary placeholder character so the insertion point is not the
# very last character slot of the paragraph (zero-length bookmarks placed
# exactly there anchor incorrectly), then add the bookmark and remove the
# placeholder while keeping the now-zero-length bookmark in place.
$rb.InsertAfter("X")
$posMark = $rb.Start
$bmRange = $d.Range($posMark, $posMark)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$placeholder = $d.Range($posMark, $posMark + 1)
$placeholder.Delete()

# ---------------------------------------------------------------------------
# Change 3: insert the graph-isomorphism parenthetical before "So in
# response..."
# ---------------------------------------------------------------------------

$r4 = $d.Content
$r4.Find.Execute("graph match the cue with the episode. So in response", $true, $false, $false, $false, $false, $true, 1, $false, "graph match the cue with the episode (i.e. determine if there exists an isomorphism between the cue and the episode). So in response", 2) | Out-Null
